$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: FAN0, FAN1, HE0, HE1, HE2 - TERM BLK 2POS SIDE ENT 3.5MM PCB
$ws.Range("D33").Value = "'"
$ws.Range("E33").Value = "'TE Connectivity"
$ws.Range("F33").Value = "'284391/2"
$ws.Range("H33").Value = "'A98159-ND"

# Row 34: HB-OUT, HB-PWR, MAIN-PWR - TERM BLK 2POS SIDE ENTRY 5MM PCB
$ws.Range("D34").Value = "'"
$ws.Range("E34").Value = "'TE Connectivity"
$ws.Range("F34").Value = "'282856-2"
$ws.Range("H34").Value = "'A98355-ND"

# Row 35: E0_MOT, E1_MOT, E2_MOT, X_MOT, Y_MOT, Z_MOT - TERM BLK 4POS SIDE ENT 3.5MM PCB
$ws.Range("D35").Value = "'"
$ws.Range("E35").Value = "'TE Connectivity"
$ws.Range("F35").Value = "'284391-4"
$ws.Range("H35").Value = "'A98161-ND"

$wb.Save()
